# Apply targeted paragraph edits (run-splitting + w:proofErr marks + lastRenderedPageBreak move)
# as produced by Word's grammar/spell checker, plus new content for the LTO question.
$d = $word.ActiveDocument

$xml86 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="3687A900" w14:textId="338BC82A" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="007C0951" w:rsidRPr="007C0951" w:rsidRDefault="007C0951" w:rsidP="00A644C2"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="007C0951"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">NEEDS FOR </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="007C0951"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>FINAL</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="007C0951"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> NIGHT</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(86).Range.InsertXML($xml86)

$xml76 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="619D6B14" w14:textId="6B76587E" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="007E521F" w:rsidRDefault="00BB7B58" w:rsidP="007E521F"><w:r><w:tab/></w:r></w:p>
'@
$d.Paragraphs.Item(76).Range.InsertXML($xml76)

$xml75 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="3B9BB041" w14:textId="215718F8" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="00BB7B58" w:rsidRDefault="00BB7B58" w:rsidP="007E521F"><w:r><w:lastRenderedPageBreak/><w:tab/><w:t>Need the login to track time?  Do we need a “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>clockout</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> button”?</w:t></w:r><w:r w:rsidR="007C0951"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="007C0951"><w:t>(Told Bob it would be)</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(75).Range.InsertXML($xml75)

$xml69 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="5265FF82" w14:textId="77777777" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="00257B5C" w:rsidRDefault="00257B5C" w:rsidP="007E521F"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:tab/></w:r><w:r><w:t xml:space="preserve">Can you auto select LTO for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sanwich</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> toppings?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$d.Paragraphs.Item(69).Range.InsertXML($xml69)

$xml55 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="06BD9AC7" w14:textId="06EF202A" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="007E521F" w:rsidRDefault="007E521F" w:rsidP="007E521F"><w:r><w:tab/><w:t xml:space="preserve">Turn the menubar into the same color as the </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>titlebar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>
'@
$d.Paragraphs.Item(55).Range.InsertXML($xml55)

$xml49 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="31D13E96" w14:textId="02DCDABD" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="00AD79A6" w:rsidRDefault="00AD79A6" w:rsidP="00F02248"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">If </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>no</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">, the window disappears, allowing them to try to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>enter info</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> again or skip </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>loyality</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$d.Paragraphs.Item(49).Range.InsertXML($xml49)

$xml41 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="2067E893" w14:textId="77777777" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="00F02248" w:rsidRDefault="00F02248" w:rsidP="00F02248"><w:r><w:tab/><w:t xml:space="preserve">Writes to database, records total </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>cost plus</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> potential tip</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(41).Range.InsertXML($xml41)

$xml35 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="58C06E53" w14:textId="77777777" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="00313090" w:rsidRDefault="00313090" w:rsidP="00313090"><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">  The underling code is </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>there,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> it just needs the correct coordinates.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(35).Range.InsertXML($xml35)

$xml21 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="63AB961E" w14:textId="7D982EB5" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="005F3653" w:rsidRDefault="005F3653" w:rsidP="005F3653"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Analytics based </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>off of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> historical trends.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(21).Range.InsertXML($xml21)

$xml17 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="04F4BF52" w14:textId="4E47D5EF" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="005F3653" w:rsidRDefault="0075294D" w:rsidP="005F3653"><w:r><w:tab/><w:t>Will need a “recall”</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>/ ”undo</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>” button</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(17).Range.InsertXML($xml17)

$xml16 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="30A1D150" w14:textId="4EE3B4B1" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="0075294D" w:rsidRDefault="0075294D" w:rsidP="005F3653"><w:r><w:tab/><w:t xml:space="preserve">Double tap to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>drop, or</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> click a button.  </w:t></w:r></w:p>
'@
$d.Paragraphs.Item(16).Range.InsertXML($xml16)

$xml13 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="126DECD6" w14:textId="49AEDCAE" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="0075294D" w:rsidRDefault="0075294D" w:rsidP="005F3653"><w:r><w:tab/><w:t xml:space="preserve">Initially based </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>off of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> number of orders</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(13).Range.InsertXML($xml13)

